# "Corrección de archivos en columnas y nombres de hojas"
#
# The workbook shipped with two sheets:
#   - "Hoja2" (sheetId 2, rId1 -> sheet1.xml) which actually holds the
#     imported Granos Básicos query-table data (A1:N807).
#   - "Hoja1" (sheetId 1, rId2 -> sheet2.xml) which is a leftover, completely
#     empty worksheet.
#
# The fix removes the stray empty sheet and renames the real data sheet from
# "Hoja2" to "Hoja1" (keeping its sheetId/relationship), so the workbook ends
# up with a single sheet named "Hoja1" holding the data, and the
# "DatosExternos_1" defined name is updated to point at the renamed sheet.

$wb = $excel.ActiveWorkbook

# Avoid any interactive "delete sheet" confirmation prompt.
$excel.DisplayAlerts = $false

# Drop the empty, superfluous worksheet.
$emptySheet = $wb.Worksheets.Item("Hoja1")
$emptySheet.Delete()

# Rename the worksheet that actually contains the data to take its place.
$dataSheet = $wb.Worksheets.Item("Hoja2")
$dataSheet.Name = "Hoja1"

$excel.DisplayAlerts = $true
